$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.830.08'
$ws.Range('E2').Value = '  -0.37%  '
$ws.Range('D3').Value = '1.639.92'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.76'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('E6').Value = '  -0.75%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.003'
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2582'
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06425'
$ws.Range('E9').Value = '  +1.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.44'
$ws.Range('E10').Value = '  +4.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07803'
$ws.Range('E11').Value = '  +0.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.271'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.644.84'
$ws.Range('E13').Value = '  +0.65%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '1.868.00'
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5620'
$ws.Range('E15').Value = '  +2.05%  '
$ws.Range('D16').Value = '0.0₅7670'
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.27'
$ws.Range('E17').Value = '  -1.10%  '
$ws.Range('D18').Value = '25.862.36'
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.004'
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.04'
$ws.Range('E20').Value = '  -1.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.379'
$ws.Range('E21').Value = '  -0.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.924'
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.142'
$ws.Range('E23').Value = '  +1.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.804'
$ws.Range('E25').Value = '  -5.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.35'
$ws.Range('E26').Value = '  -0.61%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1240'
$ws.Range('E27').Value = '  -1.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.806'
$ws.Range('E28').Value = '  +0.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.57'
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.244'
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04952'
$ws.Range('E31').Value = '  +0.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.297'
$ws.Range('E32').Value = '  +1.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.232'
$ws.Range('E33').Value = '  +1.00%  '
$ws.Range('E34').Value = '  +2.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.384'
$ws.Range('E35').Value = '  +0.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9051'
$ws.Range('E36').Value = '  +0.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5555'
$ws.Range('E37').Value = '  +0.44%  '
$ws.Range('D38').Value = '1.133.78'
$ws.Range('E38').Value = '  +1.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.551'
$ws.Range('E39').Value = '  +0.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01566'
$ws.Range('E40').Value = '  +0.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.002'
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.501'
$ws.Range('E42').Value = '  -2.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8056'
$ws.Range('E43').Value = '  +1.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '98.79'
$ws.Range('E44').Value = '  +1.50%  '
$ws.Range('D45').Value = '1.779.98'
$ws.Range('E45').Value = '  +0.62%  '
$ws.Range('D46').Value = '0.0₈111'
$ws.Range('E46').Value = '  -5.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.72'
$ws.Range('E47').Value = '  +1.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4286'
$ws.Range('E48').Value = '  -3.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.802'
$ws.Range('E49').Value = '  +3.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05040'
$ws.Range('E50').Value = '  -1.84%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9969'
$ws.Range('E51').Value = '  -0.65%  '
